# Update the 合肥-漫展信息 workbook: the scraper re-ran and the event
# list moved on by one entry (the 02.13 AG event is now in the past and
# was dropped; 02.17 becomes the new first row). All rows shift up with
# freshly refreshed stats, and the formerly-last row (05.18) is removed.
# "展览" and "全部类型" hold identical data tables and both need the edit.

$wb = $excel.ActiveWorkbook

# Row data for rows 2..8 after the refresh: B,C,D,E,F,G,H,I
# (F/G are numeric "想去人数"/"最低票价"; everything else is text)
$newRows = @(
        @("2024.02.17", "合肥·2024运动新春动漫庆典（全ip）", "锦绣大道与清潭路交口东北角 李宁体育公园", "2024.02.17 09:00-02.17 17:00", 1874, 65, "https://show.bilibili.com/platform/detail.html?id=79918", "//i0.hdslb.com/bfs/openplatform/202312/vzuMc0sJ1702902061660.jpeg"),
        @("2024.02.19", "合肥·安徽马娘only", "桐城路与庐江路交叉口西南80米 赤阑桥文玩大厦", "2024.02.19 09:00-02.19 17:00", 354, 68, "https://show.bilibili.com/platform/detail.html?id=78286", "//i1.hdslb.com/bfs/openplatform/202311/721L5pIZ1699428443216.jpeg"),
        @("2024.03.02", "合肥·星芒1.5动漫嘉年华", "山西路与太原路交叉口 挥动体育", "2024.03.02 09:30-03.02 17:30", 1145, 49.5, "https://show.bilibili.com/platform/detail.html?id=81267", "//i0.hdslb.com/bfs/openplatform/202401/GWidiefU1706003134747.jpeg"),
        @("2024.03.16", "合肥·CW国潮动漫游戏嘉年华", "南京路与庐州大道交汇处 合肥滨湖国际会展中心", "2024.03.16 09:30-03.17 17:00", 1148, 65, "https://show.bilibili.com/platform/detail.html?id=81284", "//i0.hdslb.com/bfs/openplatform/202401/38B92fWF1705995243803.jpeg"),
        @("2024.03.23", "合肥·原&铁&崩 only展", "金寨路与天堂窄路交叉口 梵木艺术中心", "2024.03.23 09:00-03.23 17:00", 50, 58, "https://show.bilibili.com/platform/detail.html?id=81574", "//i2.hdslb.com/bfs/openplatform/202401/0V5uyX6C1706697212904.png"),
        @("2024.04.04", "合肥· 第二届漫画城市动漫展 -故事再次开始", "凤淮路与固镇路西北角 庐阳全民健身中心", "2024.04.04 09:00-04.05 17:00", 5969, 60, "https://show.bilibili.com/platform/detail.html?id=78898", "//i2.hdslb.com/bfs/openplatform/202402/QnupNcrS1707125949328.jpeg"),
        @("2024.05.18", "合肥·梦时空SPO1动漫展", "阜阳路16号 银瑞林国际大酒店", "2024.05.18 10:00-05.18 17:00", 98, 60, "https://show.bilibili.com/platform/detail.html?id=80207", "//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg")
    )

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Drop the old last data row (was row 9) so the table is back to 8 rows.
    $ws.Rows.Item(9).Delete()

    for ($i = 0; $i -lt $newRows.Count; $i++) {
        $r = $i + 2
        $row = $newRows[$i]

        # Column B holds a plain-text date like "2024.02.17"; format the
        # cell as Text first so Excel doesn't coerce it into a date serial,
        # then drop back to the default "Normal" style/number format so we
        # don't leave a stray text-format style on the cell.
        $ws.Range("B$r").NumberFormat = "@"
        $ws.Range("B$r").Value = $row[0]
        $ws.Range("B$r").Style = "Normal"

        $ws.Range("C$r").Value = $row[1]
        $ws.Range("D$r").Value = $row[2]
        $ws.Range("E$r").Value = $row[3]
        $ws.Range("F$r").Value = $row[4]
        $ws.Range("G$r").Value = $row[5]
        $ws.Range("H$r").Value = $row[6]
        $ws.Range("I$r").Value = $row[7]
    }
}
